# Fruta / hortaliza, semanal
# Update the weekly Perejil (parsley) price-report rows: each data row (2-32)
# gets its Fecha (D), and for the volume/price columns (J Volumen, K Precio
# minimo, L Precio maximo, M Precio promedio ponderado, P Precio $/Kg)
# refreshed to the latest weekly figures. Some rows also swap Calidad (I)
# between "Primera" and "Segunda" as the underlying records were reordered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row=2;  D=44859; I=$null;      J=300; K=700; L=800; M=750; P=750 },
    @{ Row=3;  D=44859; I=$null;      J=200; K=600; L=600; M=600; P=600 },
    @{ Row=4;  D=44804; I=$null;      J=200; K=750; L=850; M=800; P=800 },
    @{ Row=5;  D=44804; I=$null;      J=200; K=650; L=650; M=650; P=650 },
    @{ Row=6;  D=44764; I=$null;      J=200; K=700; L=800; M=750; P=750 },
    @{ Row=7;  D=44764; I=$null;      J=150; K=600; L=600; M=600; P=600 },
    @{ Row=8;  D=44811; I=$null;      J=300; K=750; L=850; M=800; P=800 },
    @{ Row=9;  D=44868; I="Primera";  J=200; K=700; L=800; M=750; P=750 },
    @{ Row=10; D=44610; I=$null;      J=100; K=600; L=650; M=625; P=625 },
    @{ Row=11; D=44837; I=$null;      J=200; K=700; L=800; M=750; P=750 },
    @{ Row=12; D=44837; I="Segunda";  J=150; K=600; L=600; M=600; P=600 },
    @{ Row=13; D=44624; I="Primera";  J=120; K=650; L=700; M=675; P=675 },
    @{ Row=14; D=44754; I=$null;      J=200; K=700; L=750; M=725; P=725 },
    @{ Row=15; D=44608; I=$null;      J=120; K=600; L=650; M=625; P=625 },
    @{ Row=16; D=44799; I="Primera";  J=160; K=750; L=850; M=800; P=800 },
    @{ Row=17; D=44799; I="Segunda";  J=120; K=650; L=650; M=650; P=650 },
    @{ Row=18; D=44791; I="Primera";  J=240; K=750; L=800; M=775; P=775 },
    @{ Row=19; D=44791; I="Segunda";  J=250; K=650; L=650; M=650; P=650 },
    @{ Row=20; D=44831; I="Primera";  J=300; K=700; L=800; M=750; P=750 },
    @{ Row=21; D=44831; I="Segunda";  J=200; K=600; L=600; M=600; P=600 },
    @{ Row=22; D=44761; I="Primera";  J=200; K=700; L=800; M=750; P=750 },
    @{ Row=23; D=44761; I="Segunda";  J=150; K=600; L=600; M=600; P=600 },
    @{ Row=24; D=44882; I="Primera";  J=400; K=700; L=800; M=750; P=750 },
    @{ Row=25; D=44882; I="Segunda";  J=300; K=600; L=600; M=600; P=600 },
    @{ Row=26; D=44797; I=$null;      J=240; K=750; L=850; M=800; P=800 },
    @{ Row=27; D=44797; I=$null;      J=200; K=650; L=650; M=650; P=650 },
    @{ Row=28; D=44839; I=$null;      J=240; K=700; L=800; M=750; P=750 },
    @{ Row=29; D=44839; I=$null;      J=200; K=600; L=600; M=600; P=600 },
    @{ Row=30; D=44818; I=$null;      J=300; K=800; L=900; M=850; P=850 },
    @{ Row=31; D=44883; I=$null;      J=300; K=700; L=800; M=750; P=750 },
    @{ Row=32; D=44883; I="Segunda";  J=200; K=600; L=600; M=600; P=600 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 4).Value = $item.D
    if ($item.I -ne $null) {
        $ws.Cells.Item($r, 9).Value = $item.I
    }
    $ws.Cells.Item($r, 10).Value = $item.J
    $ws.Cells.Item($r, 11).Value = $item.K
    $ws.Cells.Item($r, 12).Value = $item.L
    $ws.Cells.Item($r, 13).Value = $item.M
    $ws.Cells.Item($r, 16).Value = $item.P
}
